# Auto-generated edit script applying the Zodiark_Profits sheet updates
# (scheduled runner refresh of currentAveragePrice / Leve profit calculations)
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 90910070
$ws.Range("I33").Value = 125.666664
$ws.Range("J33").Value = 200002000
$ws.Range("K33").Value = 125.666664
$ws.Range("L33").Value = 200002000
$ws.Range("M33").Value = 103.333336
$ws.Range("N33").Value = -200002458
# Row 40
$ws.Range("H40").Value = 3971
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 52
$ws.Range("H52").Value = 6439.75
$ws.Range("I52").Value = 5379.5
$ws.Range("J52").Value = 7500
$ws.Range("K52").Value = 16138.5
$ws.Range("L52").Value = 22500
$ws.Range("M52").Value = -15978.5
$ws.Range("N52").Value = -22820
# Row 68
$ws.Range("H68").Value = 458505
$ws.Range("J68").Value = 187757.5
$ws.Range("L68").Value = 187757.5
$ws.Range("N68").Value = -189255.5
# Row 71
$ws.Range("H71").Value = 458505
$ws.Range("J71").Value = 187757.5
$ws.Range("L71").Value = 563272.5
$ws.Range("N71").Value = -570760.5
# Row 138
$ws.Range("H138").Value = 1841.585
$ws.Range("I138").Value = 1053.7894
$ws.Range("J138").Value = 2281.8235
$ws.Range("K138").Value = 3161.3682
$ws.Range("L138").Value = 6845.470499999999
$ws.Range("M138").Value = 1978.6318
$ws.Range("N138").Value = -17125.4705

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1469.5151
$ws.Range("I2").Value = 715.1111
$ws.Range("K2").Value = 715.1111
$ws.Range("M2").Value = -602.1111
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 116
$ws.Range("H116").Value = 1469.5151
$ws.Range("I116").Value = 715.1111
$ws.Range("K116").Value = 715.1111
$ws.Range("M116").Value = 1578.8889

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1469.5151
$ws.Range("I3").Value = 715.1111
$ws.Range("K3").Value = 715.1111
$ws.Range("M3").Value = -601.1111
# Row 15
$ws.Range("H15").Value = 33337
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 99999
$ws.Range("K15").Value = 6
$ws.Range("L15").Value = 99999
$ws.Range("M15").Value = 221
$ws.Range("N15").Value = -100453
# Row 134
$ws.Range("H134").Value = 3426.842
$ws.Range("I134").Value = 3411.7222
$ws.Range("K134").Value = 10235.1666
$ws.Range("M134").Value = -7700.1666

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 602.8929000000001
$ws.Range("I22").Value = 245.75
$ws.Range("J22").Value = 1079.0834
$ws.Range("K22").Value = 245.75
$ws.Range("L22").Value = 1079.0834
$ws.Range("M22").Value = 104.25
$ws.Range("N22").Value = -1779.0834
# Row 31
$ws.Range("H31").Value = 2558.3333
$ws.Range("I31").Value = 2675.7856
$ws.Range("J31").Value = 914
$ws.Range("K31").Value = 2675.7856
$ws.Range("L31").Value = 914
$ws.Range("M31").Value = -2380.7856
$ws.Range("N31").Value = -1504
# Row 34
$ws.Range("H34").Value = 2558.3333
$ws.Range("I34").Value = 2675.7856
$ws.Range("J34").Value = 914
$ws.Range("K34").Value = 2675.7856
$ws.Range("L34").Value = 914
$ws.Range("M34").Value = -2473.7856
$ws.Range("N34").Value = -1318
# Row 98
$ws.Range("H98").Value = 22332.666
$ws.Range("J98").Value = 22332.666
$ws.Range("L98").Value = 22332.666
$ws.Range("N98").Value = -26824.666
# Row 99
$ws.Range("H99").Value = 3273.5454
$ws.Range("I99").Value = 3156
$ws.Range("J99").Value = 3479.25
$ws.Range("K99").Value = 3156
$ws.Range("L99").Value = 3479.25
$ws.Range("M99").Value = -1658
$ws.Range("N99").Value = -6475.25
# Row 126
$ws.Range("H126").Value = 3273.5454
$ws.Range("I126").Value = 3156
$ws.Range("J126").Value = 3479.25
$ws.Range("K126").Value = 9468
$ws.Range("L126").Value = 10437.75
$ws.Range("M126").Value = -6998
$ws.Range("N126").Value = -15377.75
# Row 134
$ws.Range("H134").Value = 2804
$ws.Range("J134").Value = 3002.2
$ws.Range("L134").Value = 9006.599999999999
$ws.Range("N134").Value = -14076.6

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 861.36365
$ws.Range("J22").Value = 875
$ws.Range("L22").Value = 2625
$ws.Range("N22").Value = -2963
# Row 27
$ws.Range("H27").Value = 861.36365
$ws.Range("J27").Value = 875
$ws.Range("L27").Value = 2625
$ws.Range("N27").Value = -2829
# Row 35
$ws.Range("H35").Value = 466.33334
$ws.Range("J35").Value = 698.5
$ws.Range("L35").Value = 2095.5
$ws.Range("N35").Value = -2671.5
# Row 41
$ws.Range("H41").Value = 3467
$ws.Range("J41").Value = 3701
$ws.Range("L41").Value = 11103
$ws.Range("N41").Value = -11779
# Row 62
$ws.Range("H62").Value = 4979.5
$ws.Range("J62").Value = 8999
$ws.Range("L62").Value = 26997
$ws.Range("N62").Value = -28369
# Row 65
$ws.Range("H65").Value = 4979.5
$ws.Range("J65").Value = 8999
$ws.Range("L65").Value = 80991
$ws.Range("N65").Value = -87855
# Row 69
$ws.Range("H69").Value = 1400
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1400
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 4200
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -5822
# Row 72
$ws.Range("H72").Value = 1400
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1400
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 12600
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -20712
# Row 107
$ws.Range("H107").Value = 840.6842
$ws.Range("I107").Value = 641.6667
$ws.Range("K107").Value = 1925.0001
$ws.Range("M107").Value = -5.000099999999975

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3306.2083
$ws.Range("I102").Value = 3064.3333
$ws.Range("J102").Value = 4999.3335
$ws.Range("K102").Value = 3064.3333
$ws.Range("L102").Value = 4999.3335
$ws.Range("M102").Value = -1442.3333
$ws.Range("N102").Value = -8243.333500000001
# Row 126
$ws.Range("H126").Value = 3203.2222
$ws.Range("J126").Value = 5116.6665
$ws.Range("L126").Value = 15349.9995
$ws.Range("N126").Value = -20289.9995

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 109
$ws.Range("H109").Value = 36750
$ws.Range("J109").Value = 36750
$ws.Range("L109").Value = 36750
$ws.Range("N109").Value = -39524

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -888
# Row 37
$ws.Range("H37").Value = 54350
$ws.Range("I37").Value = 48700
$ws.Range("J37").Value = 60000
$ws.Range("K37").Value = 48700
$ws.Range("L37").Value = 60000
$ws.Range("M37").Value = -48497
$ws.Range("N37").Value = -60406
# Row 42
$ws.Range("H42").Value = 18999.5
$ws.Range("J42").Value = 18999.5
$ws.Range("L42").Value = 18999.5
$ws.Range("N42").Value = -19755.5
# Row 43
$ws.Range("H43").Value = 34524
$ws.Range("I43").Value = 32699
$ws.Range("K43").Value = 32699
$ws.Range("M43").Value = -32550
